$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("CL001", "Johnson", 0.33, 13.33),
    @("CL007", "Shah/Rieke", 1.33, 53.33),
    @("CL010", "Silverton", 0.33, 13.33)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}
